# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6013.6
$ws.Range("I19").Value = 9269.200000000001
$ws.Range("K19").Value = 9269.200000000001
$ws.Range("M19").Value = -9094.200000000001

$ws.Range("H33").Value = 2651.2144
$ws.Range("I33").Value = 1720.1818
$ws.Range("K33").Value = 1720.1818
$ws.Range("M33").Value = -1491.1818

$ws.Range("H51").Value = 6947842
$ws.Range("I51").Value = 3692.5
$ws.Range("J51").Value = 10419917
$ws.Range("K51").Value = 3692.5
$ws.Range("L51").Value = 10419917
$ws.Range("M51").Value = -3208.5
$ws.Range("N51").Value = -10420885

$ws.Range("H86").Value = 2707
$ws.Range("J86").Value = 2611.3333
$ws.Range("L86").Value = 2611.3333
$ws.Range("N86").Value = -4857.3333

$ws.Range("H89").Value = 2707
$ws.Range("J89").Value = 2611.3333
$ws.Range("L89").Value = 13056.6665
$ws.Range("N89").Value = -24288.6665

$ws.Range("H112").Value = 2011.3334
$ws.Range("J112").Value = 2273.182
$ws.Range("L112").Value = 6819.545999999999
$ws.Range("N112").Value = -9035.545999999998

$ws.Range("H134").Value = 105308.27
$ws.Range("J134").Value = 97889.3
$ws.Range("L134").Value = 97889.3
$ws.Range("N134").Value = -108029.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 717899.8
$ws.Range("I45").Value = 1670417.9
$ws.Range("K45").Value = 1670417.9
$ws.Range("M45").Value = -1670040.9

$ws.Range("H132").Value = 10742.591
$ws.Range("I132").Value = 12842.765
$ws.Range("K132").Value = 38528.295
$ws.Range("M132").Value = -35998.295

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 96686
$ws.Range("J58").Value = 96686
$ws.Range("L58").Value = 96686
$ws.Range("N58").Value = -97274

$ws.Range("H86").Value = 2179.4614
$ws.Range("I86").Value = 1975.7333
$ws.Range("K86").Value = 1975.7333
$ws.Range("M86").Value = -852.7333000000001

$ws.Range("H89").Value = 2179.4614
$ws.Range("I89").Value = 1975.7333
$ws.Range("K89").Value = 9878.666500000001
$ws.Range("M89").Value = -4262.666500000001

$ws.Range("H105").Value = 4352554.5
$ws.Range("I105").Value = 4766749
$ws.Range("K105").Value = 4766749
$ws.Range("M105").Value = -4765002

$ws.Range("H134").Value = 2030.5238
$ws.Range("I134").Value = 1545.75
$ws.Range("J134").Value = 3581.8
$ws.Range("K134").Value = 4637.25
$ws.Range("L134").Value = 10745.4
$ws.Range("M134").Value = -2102.25
$ws.Range("N134").Value = -15815.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 260
$ws.Range("I22").Value = 260
$ws.Range("K22").Value = 260
$ws.Range("M22").Value = 90

$ws.Range("H62").Value = 12471
$ws.Range("I62").Value = 13418.1
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 13418.1
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -12794.1
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 12471
$ws.Range("I65").Value = 13418.1
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 67090.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -63970.5
$ws.Range("N65").Value = -21240

$ws.Range("H69").Value = 16333.333
$ws.Range("I69").Value = 9000
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 20000
$ws.Range("M69").Value = -8251
$ws.Range("N69").Value = -21498

$ws.Range("H72").Value = 16333.333
$ws.Range("I72").Value = 9000
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 60000
$ws.Range("M72").Value = -23256
$ws.Range("N72").Value = -67488

$ws.Range("H99").Value = 20456
$ws.Range("J99").Value = 2199
$ws.Range("L99").Value = 2199
$ws.Range("N99").Value = -5195

$ws.Range("H105").Value = 898.75
$ws.Range("I105").Value = 748.5714
$ws.Range("K105").Value = 748.5714
$ws.Range("M105").Value = 998.4286

$ws.Range("H126").Value = 20456
$ws.Range("J126").Value = 2199
$ws.Range("L126").Value = 6597
$ws.Range("N126").Value = -11537

$ws.Range("H130").Value = 100780
$ws.Range("J130").Value = 100780
$ws.Range("L130").Value = 100780
$ws.Range("N130").Value = -110820

$ws.Range("H132").Value = 2069.2083
$ws.Range("I132").Value = 2121.8635
$ws.Range("K132").Value = 6365.5905
$ws.Range("M132").Value = -3835.5905

$ws.Range("H134").Value = 4394.3
$ws.Range("I134").Value = 4743
$ws.Range("K134").Value = 14229
$ws.Range("M134").Value = -11694

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1047.8948
$ws.Range("J122").Value = 1093.871
$ws.Range("L122").Value = 9844.839
$ws.Range("N122").Value = -14744.839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 277.26923
$ws.Range("I2").Value = 237.47058
$ws.Range("J2").Value = 352.44446
$ws.Range("K2").Value = 237.47058
$ws.Range("L2").Value = 352.44446
$ws.Range("M2").Value = -124.47058
$ws.Range("N2").Value = -578.4444599999999

$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("N33").Value = -25504

$ws.Range("H70").Value = 6500.4443
$ws.Range("I70").Value = 6267.5
$ws.Range("J70").Value = 6966.3335
$ws.Range("K70").Value = 6267.5
$ws.Range("L70").Value = 6966.3335
$ws.Range("M70").Value = -5997.5
$ws.Range("N70").Value = -7506.3335

$ws.Range("H73").Value = 6500.4443
$ws.Range("I73").Value = 6267.5
$ws.Range("J73").Value = 6966.3335
$ws.Range("K73").Value = 6267.5
$ws.Range("L73").Value = 6966.3335
$ws.Range("M73").Value = -5331.5
$ws.Range("N73").Value = -8838.333500000001

$ws.Range("H80").Value = 4006.7144
$ws.Range("I80").Value = 3500.5715
$ws.Range("J80").Value = 4512.857
$ws.Range("K80").Value = 3500.5715
$ws.Range("L80").Value = 4512.857
$ws.Range("M80").Value = -2502.5715
$ws.Range("N80").Value = -6508.857

$ws.Range("H83").Value = 4006.7144
$ws.Range("I83").Value = 3500.5715
$ws.Range("J83").Value = 4512.857
$ws.Range("K83").Value = 17502.8575
$ws.Range("L83").Value = 22564.285
$ws.Range("M83").Value = -12510.8575
$ws.Range("N83").Value = -32548.285

$ws.Range("H107").Value = 33334832
$ws.Range("J107").Value = 50001748
$ws.Range("L107").Value = 50001748
$ws.Range("N107").Value = -50005588

$ws.Range("H132").Value = 2409.5386
$ws.Range("I132").Value = 2467.3914
$ws.Range("K132").Value = 7402.174199999999
$ws.Range("M132").Value = -4872.174199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5224.75
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5224.75
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15674.25
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -20574.25

$ws.Range("H133").Value = 87225.25
$ws.Range("J133").Value = 87225.25
$ws.Range("L133").Value = 87225.25
$ws.Range("N133").Value = -92285.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30300.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 30300.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 30300.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -32296.5

$ws.Range("H83").Value = 30300.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 30300.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 90901.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -100885.5

$ws.Range("H120").Value = 8000
$ws.Range("J120").Value = 8000
$ws.Range("L120").Value = 8000
$ws.Range("N120").Value = -17676

$ws.Range("H126").Value = 3658.4
$ws.Range("I126").Value = 3521.4
$ws.Range("K126").Value = 10564.2
$ws.Range("M126").Value = -8094.200000000001
